# Lost Creek Past Reports.docx — add proofing marks (spell/grammar check
# artifacts) around specific words/phrases, and append a new fishing
# report at the end of the document.
#
# Strategy: Word's COM Range.InsertXML() call replaces the *entire*
# paragraph(s) addressed by the Range when given a <w:p>-rooted OOXML
# fragment (bare, non-<w:p> fragments corrupt the paragraph, so every
# fragment below is always wrapped in one or more full <w:p> elements).
# So: locate the paragraph containing the target sentence, then rebuild
# that whole paragraph's XML with the sentence split into extra runs
# around <w:proofErr/> markers, preserving any other runs untouched.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParagraphByText {
    param($d, [string]$needle)
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphXml {
    param($p, [string]$innerXml)
    $rng = $p.Range
    $xml = "<w:p $wNs>" + $innerXml + "</w:p>"
    $rng.InsertXML($xml)
}

# NOTE: this COM-interop PowerShell parser mis-binds a parenthesized
# expression `(...)` given as the 2nd positional argument when the 1st
# argument is a COM object (e.g. `Set-ParagraphXml $p (...)` silently
# yields $null for both params). Always build the inner-XML string into
# a plain variable first, then pass that variable - never an inline
# parenthesized expression - as the argument.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "red and gold Pot-o-Golds, red and white Kamloopers, and silver
#    Crocodile spoons work"  -> spellStart/spellEnd around "Kamloopers"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "red and gold Pot-o-Golds, red and white Kamloopers, and silver Crocodile spoons work"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">red and gold Pot-o-Golds, red and white </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamloopers</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, and silver Crocodile spoons work</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 2) "of the fish are planter sized Rainbows and Cutts. We also did
#    catch a Tiger Trout." -> spellStart/spellEnd around "Cutts"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "of the fish are planter sized Rainbows and Cutts. We also did catch a Tiger Trout."
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">of the fish are planter sized Rainbows and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cutts</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>. We also did catch a Tiger Trout.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 3) "any fish is to use Kamlooper spoons with the best colors being red
#    and gold or" -> spellStart/spellEnd around "Kamlooper"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "any fish is to use Kamlooper spoons with the best colors being red and gold or"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">any fish is to use </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamlooper</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> spoons with the best colors being red and gold or</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 4) "will catch you a good amount of fish. While it is clam find an
#    area with a far" -> gramStart/gramEnd around "amount"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "will catch you a good amount of fish. While it is clam find an area with a far"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">will catch you a good </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>amount</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> of fish. While it is clam find an area with a far</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 5) "small spoons like junior Kamloopers and Pot-o-Golds will also get
#    some bites." -> spellStart/spellEnd around "Kamloopers"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "small spoons like junior Kamloopers and Pot-o-Golds will also get some bites."
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">small spoons like junior </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamloopers</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and Pot-o-Golds will also get some bites.</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 6) "marshmallow on the bottom with a 2 to 3 foot leader. You won't
#    need to cast your" -> gramStart/gramEnd around "2 to 3 foot"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "marshmallow on the bottom with a 2 to 3 foot leader. You won"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">marshmallow on the bottom with a </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>2 to 3 foot</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> leader. You won' + [char]0x2019 + 't need to cast your</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 7) "green Jake, and a follow on a Kamlooper. Beyond that, there was
#    absolutely no action. " -> spellStart/spellEnd around "Kamlooper"
#    (this paragraph also carries a trailing bold "(Report from: ...)"
#    run plus a bold paragraph-mark rPr, both preserved as-is)
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "green Jake, and a follow on a Kamlooper. Beyond that, there was absolutely no action. "
Set-ParagraphXml $p (
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">green Jake, and a follow on a </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamlooper</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">. Beyond that, there was absolutely no action. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Report from: 4/4/20)</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 8) "worm and half or whole marshmallow. Spoons like Kamloopers would
#    also work just" -> spellStart/spellEnd around "Kamloopers"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "worm and half or whole marshmallow. Spoons like Kamloopers would also work just"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">worm and half or whole marshmallow. Spoons like </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamloopers</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> would also work just</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 9) "marshmallow on about a 1 to 3 foot leader works best. " ->
#    gramStart/gramEnd around "1 to 3 foot"
#    (this paragraph also carries a trailing bold "(Report From: ...)"
#    run plus a bold paragraph-mark rPr, both preserved as-is)
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "marshmallow on about a 1 to 3 foot leader works best. "
Set-ParagraphXml $p (
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">marshmallow on about a </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>1 to 3 foot</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> leader works best. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Report From: 4/18/21)</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 10) "spoons for us were either yellow and gold or red white and gold
#     Kamloopers. And" -> spellStart/spellEnd around "Kamloopers"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "spoons for us were either yellow and gold or red white and gold Kamloopers. And"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">spoons for us were either yellow and gold or red white and gold </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kamloopers</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>. And</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 11) "jigs and small yellow Atomic tubes about 1 to 2' below the
#     bobber tipped with" -> gramStart/gramEnd around "Atomic"
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "jigs and small yellow Atomic tubes about 1 to 2"
Set-ParagraphXml $p (
    '<w:r><w:t xml:space="preserve">jigs and small yellow </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Atomic</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> tubes about 1 to 2' + [char]0x27 + ' below the bobber tipped with</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 12) "pink maggots or salmon eggs. (Report From: 11/7/21)" paragraph
#     gains a bold paragraph-mark rPr (no text changes in this
#     paragraph itself).
# ---------------------------------------------------------------------
$p = Get-ParagraphByText $d "pink maggots or salmon eggs."
Set-ParagraphXml $p (
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">pink maggots or salmon eggs. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Report From: 11/7/21)</w:t></w:r>'
)

# ---------------------------------------------------------------------
# 13) Append a brand-new fishing report at the end of the document,
#     after a blank bold-marked separator paragraph, replacing the
#     final (empty) paragraph and re-adding an empty one after it so
#     the document still ends on a blank paragraph.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$tailRng = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$newParasXml =
    "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>" +
    "<w:p $wNs><w:r><w:t>Fishing is beginning to slow down a little probably due to the fact that it</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>is getting late in the season. The best method right now is casting a 2.5" + [char]0x22 + " white</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>tube jig tipped with Chub meat. The other method that works is fishing on the</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>bottom with worm and a marshmallow on a 1 to 3" + [char]0x27 + " leader. Most of the fish are close</w:t></w:r></w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:t>to the shore, so you do not have to cast out very far at all to catch the fish.</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Report From: 11/2</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1</w:t></w:r>" +
        "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>/21)</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs/>"

$tailRng.InsertXML($newParasXml)
